$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F36").Value = 78
$ws.Range("G36").Value = 15348.06

$ws.Range("F41").Value = 202
$ws.Range("G41").Value = 38963.78

$ws.Range("F48").Value = 229
$ws.Range("G48").Value = 12846.9

$ws.Range("F51").Value = 134
$ws.Range("G51").Value = 12534.36

$ws.Range("F55").Value = 114
$ws.Range("G55").Value = 6356.64

$ws.Range("F56").Value = 38
$ws.Range("G56").Value = 848.16

$ws.Range("F61").Value = 206
$ws.Range("G61").Value = 53710.38

$ws.Range("F64").Value = 58
$ws.Range("G64").Value = 4612.16

$ws.Range("B66").Value = 189260.89

$ws.Range("B126").Value = 65258
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 0

$ws.Range("B127").Value = 64196
$ws.Range("F127").Value = 1
$ws.Range("G127").Value = 32143.58

$ws.Range("F141").Value = 41
$ws.Range("G141").Value = 2134.46

$ws.Range("B147").Value = 20488.83

$ws.Range("B161").Value = 57756
$ws.Range("E161").Value = 79.37
$ws.Range("F161").Value = -100
$ws.Range("G161").Value = -6644

$ws.Range("B162").Value = 64350
$ws.Range("E162").Value = 70.63
$ws.Range("F162").Value = 2
$ws.Range("G162").Value = 132.88

$ws.Range("F195").Value = 8
$ws.Range("G195").Value = 600.16

$ws.Range("B197").Value = 1005.92

$ws.Range("F210").Value = 101
$ws.Range("G210").Value = 5495.41

$ws.Range("B218").Value = 71986.03

$ws.Range("F222").Value = 606
$ws.Range("G222").Value = 11211

$ws.Range("B229").Value = 20392.4

$ws.Range("F262").Value = 69
$ws.Range("G262").Value = 5404.77

$ws.Range("F268").Value = 8
$ws.Range("G268").Value = 1017.68

$ws.Range("F280").Value = 12
$ws.Range("G280").Value = 1165.2

$ws.Range("F284").Value = 145
$ws.Range("G284").Value = 6796.15

$ws.Range("F287").Value = 45
$ws.Range("G287").Value = 2463.3

$ws.Range("B290").Value = 66194
$ws.Range("C290").Value = "HIM-Total Care Baby Pants Diapers-M-9s"
$ws.Range("F290").Value = 22
$ws.Range("G290").Value = 1884.96

$ws.Range("B291").Value = 64983
$ws.Range("C291").Value = "HIM-TOTAL CARE BABY PANTS DIAPERS-M-9S"
$ws.Range("F291").Value = 6
$ws.Range("G291").Value = 514.08

$ws.Range("B295").Value = 108170.85

$ws.Range("B306").Value = 63531
$ws.Range("E306").Value = 152.53
$ws.Range("F306").Value = 26
$ws.Range("G306").Value = 3730.48

$ws.Range("B307").Value = 57802
$ws.Range("E307").Value = 162.71
$ws.Range("F307").Value = -79
$ws.Range("G307").Value = -11334.92

$ws.Range("B308").Value = 63510
$ws.Range("E308").Value = 50.66
$ws.Range("F308").Value = 76
$ws.Range("G308").Value = 3620.64

$ws.Range("B309").Value = 55356
$ws.Range("E309").Value = 54.04
$ws.Range("F309").Value = -158
$ws.Range("G309").Value = -7527.12

$ws.Range("B317").Value = 60325
$ws.Range("E317").Value = 151.57
$ws.Range("F317").Value = -102
$ws.Range("G317").Value = -12939.72

$ws.Range("B318").Value = 63560
$ws.Range("E318").Value = 134.87
$ws.Range("F318").Value = 1
$ws.Range("G318").Value = 126.86

$ws.Range("F324").Value = 7
$ws.Range("G324").Value = 1199.31

$ws.Range("B328").Value = -17716.82

$ws.Range("F366").Value = 49
$ws.Range("G366").Value = 2711.17

$ws.Range("F367").Value = 122
$ws.Range("G367").Value = 7405.4

$ws.Range("F370").Value = 195
$ws.Range("G370").Value = 32368.05

$ws.Range("B372").Value = 55082.22

$ws.Range("B381").Value = 58047
$ws.Range("D381").Value = 105.54
$ws.Range("E381").Value = 126.1
$ws.Range("F381").Value = 32
$ws.Range("G381").Value = 3377.28

$ws.Range("B382").Value = 47097
$ws.Range("D382").Value = 112.28
$ws.Range("E382").Value = 134.16
$ws.Range("F382").Value = 15
$ws.Range("G382").Value = 1684.2

$ws.Range("F384").Value = 4
$ws.Range("G384").Value = 1158.64

$ws.Range("F387").Value = 399
$ws.Range("G387").Value = 38543.4

$ws.Range("B389").Value = 55025.88

$ws.Range("F408").Value = 190
$ws.Range("G408").Value = 3011.5

$ws.Range("B417").Value = 163448.18

$ws.Range("F431").Value = 18
$ws.Range("G431").Value = 484.02

$ws.Range("F432").Value = 76
$ws.Range("G432").Value = 3679.16

$ws.Range("B438").Value = 23461.04

$ws.Range("F453").Value = 29
$ws.Range("G453").Value = 4214.28

$ws.Range("B458").Value = 89098.74000000001

$ws.Range("B479").Value = 64810
$ws.Range("E479").Value = 291.22
$ws.Range("F479").Value = 0
$ws.Range("G479").Value = 0

$ws.Range("B480").Value = 53319
$ws.Range("E480").Value = 310.64
$ws.Range("F480").Value = -6
$ws.Range("G480").Value = -1643.52

$ws.Range("F511").Value = 203
$ws.Range("G511").Value = 20273.61

$ws.Range("F523").Value = 139
$ws.Range("G523").Value = 11899.79

$ws.Range("B525").Value = 115127.85

$ws.Range("F527").Value = 37
$ws.Range("G527").Value = 1225.07

$ws.Range("F529").Value = 114
$ws.Range("G529").Value = 3774.54

$ws.Range("F531").Value = 205
$ws.Range("G531").Value = 6787.55

$ws.Range("B535").Value = 21795.3

$ws.Range("F555").Value = 4
$ws.Range("G555").Value = 1016.04

$ws.Range("B556").Value = 40848.71

$ws.Range("F558").Value = 166
$ws.Range("G558").Value = 20227.1

$ws.Range("B561").Value = 23682.06

$ws.Range("F606").Value = 0
$ws.Range("G606").Value = 0

$ws.Range("B607").Value = 22627

$ws.Range("F622").Value = 475
$ws.Range("G622").Value = 48882.25

$ws.Range("F625").Value = 310
$ws.Range("G625").Value = 11417.3

$ws.Range("B628").Value = 201009.74

$ws.Range("F660").Value = 45
$ws.Range("G660").Value = 1338.3

$ws.Range("B668").Value = 11009.81

$ws.Range("F704").Value = 7
$ws.Range("G704").Value = 2499.7

$ws.Range("B713").Value = 61196.77

$ws.Range("B718").Value = 2464563.96

$ws.Range("B719").Value = 2464563.96
